$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: custom_nvarchar (NOT NULL) - length shrinks from 555 to 250 (MaxLength metadata is in bytes: 250*2 = 500)
$ws.Range("G13").Value = 500
$ws.Range("J13").Value = "CREATE TYPE [TestSchema].[custom_nvarchar]`n    FROM nvarchar(250) NOT NULL`nGO`n"

# Row 11: custom_nchar - length grows from 9 to 45 (MaxLength metadata is in bytes: 45*2 = 90)
$ws.Range("G11").Value = 90
$ws.Range("J11").Value = "CREATE TYPE [TestSchema].[custom_nchar]`n    FROM nchar(45) NULL`nGO`n"

# Move the last active selection in the frozen (bottom-left) pane to J12, matching the
# cell the author was last working in.
$ws.Range("J12").Select()
